$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title in A1 (year 2023 -> 2025)
$ws.Range("A1").Value = "Peak Loads for 2025 (Mthembanji)"

# Update the peak demand values for rows 3-6
$ws.Range("B3").Value = 28.3032575571
$ws.Range("B4").Value = 26.2917684729
$ws.Range("B5").Value = 26.9020471434
$ws.Range("B6").Value = 0
